$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.556.11"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +6.00%  "
$ws.Range("D3").Value = "'3.564.07"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'413.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").Value = "'128.82"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("D7").Value = "'0.646"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +3.41%  "
$ws.Range("D8").Value = "'3.557.98"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.63%  "
$ws.Range("E9").Value = "  +0.00%  "
$ws.Range("D10").Value = "'0.768"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.96%  "
$ws.Range("D11").Value = "'0.174"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +17.31%  "
$ws.Range("D12").Value = "'0.0000322"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +46.63%  "
$ws.Range("D13").Value = "'42.07"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.23%  "
$ws.Range("D14").Value = "'9.81"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").Value = "'4.121.48"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("E16").Value = "  -0.41%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'20.13"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.81%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "'3.610.52"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.48%  "
$ws.Range("D19").Value = "'1.11"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.98%  "
$ws.Range("D20").Value = "'66.435.66"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +5.91%  "
$ws.Range("D21").Value = "'12.22"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.31%  "
$ws.Range("D22").Value = "'446.98"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.13%  "
$ws.Range("D23").Value = "'89.03"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").Value = "'3.13"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.16%  "
$ws.Range("D25").Value = "'12.94"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.74%  "
$ws.Range("D26").Value = "'3.30"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").Value = "'9.90"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -6.12%  "
$ws.Range("D28").Value = "'34.65"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.70%  "
$ws.Range("E29").Value = "  +1.24%  "
$ws.Range("D30").Value = "'2.77"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.37%  "
$ws.Range("D31").Value = "'12.27"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.29%  "
$ws.Range("D32").Value = "'0.116"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +3.49%  "
$ws.Range("D33").Value = "'7.27"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.11%  "
$ws.Range("E34").Value = "  -4.49%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'39.46"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -3.35%  "
$ws.Range("D37").Value = "'56.42"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.59%  "
$ws.Range("D38").Value = "'0.0487"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("D39").Value = "'0.0₃0733"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +31.26%  "
$ws.Range("D40").Value = "'0.148"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +10.47%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").Value = "'147.89"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.13%  "
$ws.Range("E43").Value = "  -4.13%  "
$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'2.71"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("E45").Value = "  -3.43%  "
$ws.Range("D46").Value = "'4.29"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.05%  "
$ws.Range("D47").Value = "'0.306"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -4.62%  "
$ws.Range("D48").Value = "'1.96"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -5.75%  "
$ws.Range("D49").Value = "'2.26"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.91%  "
$ws.Range("D50").Value = "'116.47"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +6.42%  "
$ws.Range("D51").Value = "'15.46"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.33%  "
